$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estudiantes")

# Delete row 3 (Juan Perez) entirely, shifting remaining rows up.
$ws.Rows("3:3").Delete()

# The former row 4 (Lujan Gomez) is now row 3; update last name to Martínez.
$ws.Range("C3").Value = "Martínez"
